# Fruta / hortaliza, semanal
#
# Inserts 3 new weekly rows (new rows 424-426) just above the existing
# row that used to be "424" (date 44249 / Primera / Diguillín), pushing
# the old rows 424-432 down to 427-435. The new rows carry a fresh
# report date (44939) for three quality grades: Especial, Primera and
# Segunda, all priced at the new "$8000 / caja 7 kilos" level (plus the
# Especial entry). Dimension grows from A1:T432 to A1:T435.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 424-432 down by inserting 3 blank rows at 424.
$ws.Range("A424:T426").Insert()

# --- New row 424: Especial ---
$ws.Range("A424").Value = 7
$ws.Range("B424").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C424").Value = "Ñuble"
$ws.Range("D424").Value = 44939
$ws.Range("E424").Value = 16
$ws.Range("F424").Value = "Fruta"
$ws.Range("G424").Value = 100101
$ws.Range("H424").Value = "Berries"
$ws.Range("I424").Value = 100112025
$ws.Range("J424").Value = "Frutilla"
$ws.Range("K424").Value = "Sin especificar"
$ws.Range("L424").Value = "Especial"
$ws.Range("M424").Value = 80
$ws.Range("N424").Value = 8000
$ws.Range("O424").Value = 8000
$ws.Range("P424").Value = 8000
$ws.Range("Q424").Value = "$/caja 7 kilos"
$ws.Range("R424").Value = "Provincia de Diguillín"
$ws.Range("S424").Value = 1143
$ws.Range("T424").Value = 7

# --- New row 425: Primera ---
$ws.Range("A425").Value = 7
$ws.Range("B425").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C425").Value = "Ñuble"
$ws.Range("D425").Value = 44939
$ws.Range("E425").Value = 16
$ws.Range("F425").Value = "Fruta"
$ws.Range("G425").Value = 100101
$ws.Range("H425").Value = "Berries"
$ws.Range("I425").Value = 100112025
$ws.Range("J425").Value = "Frutilla"
$ws.Range("K425").Value = "Sin especificar"
$ws.Range("L425").Value = "Primera"
$ws.Range("M425").Value = 80
$ws.Range("N425").Value = 7000
$ws.Range("O425").Value = 7000
$ws.Range("P425").Value = 7000
$ws.Range("Q425").Value = "$/caja 7 kilos"
$ws.Range("R425").Value = "Provincia de Diguillín"
$ws.Range("S425").Value = 1000
$ws.Range("T425").Value = 7

# --- New row 426: Segunda ---
$ws.Range("A426").Value = 7
$ws.Range("B426").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C426").Value = "Ñuble"
$ws.Range("D426").Value = 44939
$ws.Range("E426").Value = 16
$ws.Range("F426").Value = "Fruta"
$ws.Range("G426").Value = 100101
$ws.Range("H426").Value = "Berries"
$ws.Range("I426").Value = 100112025
$ws.Range("J426").Value = "Frutilla"
$ws.Range("K426").Value = "Sin especificar"
$ws.Range("L426").Value = "Segunda"
$ws.Range("M426").Value = 80
$ws.Range("N426").Value = 6000
$ws.Range("O426").Value = 6000
$ws.Range("P426").Value = 6000
$ws.Range("Q426").Value = "$/caja 7 kilos"
$ws.Range("R426").Value = "Provincia de Diguillín"
$ws.Range("S426").Value = 857
$ws.Range("T426").Value = 7
